$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the unique resource config row (row 2)
$ws.Range("B2").Value = "UAlbany Career Services"
$ws.Range("D2").Value = "UAlbany Career Services"
$ws.Range("C2").Value = "ACAdata()"
$ws.Range("A2").Value = "ACA_CareerServices"
$ws.Range("E2").Value = "T"
$ws.Range("F2").Value = "T"

# Update the static "ENTER INTO UI" text cells (A10:A12) to mirror the
# recomputed formula results in B4:B6
$ws.Range("A10").Value = "mod_Accordion_ui('ACA_CareerServices')"
$ws.Range("A11").Value = "mod_Accordion_server('ACA_CareerServices', selector=selection, data=ACAdata(), title = c('UAlbany Career Services'), Visible = T)"
$ws.Range("A12").Value = "mod_info_server('ACA_CareerServices', selector = selection, data = ACAdata(), rownametitle = c('UAlbany Career Services'), phone = T, website = T)"

# Match the new selection state recorded in the sheet view
$ws.Range("B8").Select()
